$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for rule R10 (cell E8) from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active cell/selection on the sheet (E8), matching the author's saved state
$ws.Activate()
$ws.Range("E8").Select()
